$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "⚡️🚀המאמר היומי של מייק 08.07.24: ⚡️🚀", $true, $false, $false, $false, $false,
    $true, 1, $false, "⚡️🚀המאמר היומי של מייק 07.07.24:⚡️🚀", 2)

$d.Content.Find.Execute(
    "Mixture of A Million Experts", $true, $false, $false, $false, $false,
    $true, 1, $false, "The Road Less Scheduled", 2)

$d.Content.Find.Execute(
    "המאמר של היום מציע לקחת את שיטת (Mixture of Experts(MoE לבניית ארכיטקטורות של מודלים עמוקים פופולרית במיוחד במודלי שפה. מאוד בגדול ב- MoE הרשת מורכבת מתת-רשתות (בד״כ מחלקים את שכבת ה-FFN של הטרנספורמר לכמה חלקים זרים). MoE מאומן להשתמש כל בפעם בחלק מתת-רשתות אלו (הנקראות מומחים) כאשר רשת gating רדודה יחסית באיזה מומחים צריך להשתמש כל פעם. כלומר יש לנו כן סוג של מימוש הגישה שנקראת ״lottery ticket hypothesis`" דינמי כאשר כל פעם בוחרים להריץ רק חלק מהרשת. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "היום סוקרים מאמר שלא נראה כמו מאמר למידה עמוקה רגיל. בהתחלה זה אולי יכול להיראות שהמאמר מציע עוד שכלול מי יודע מה ל-ADAM או שיטה אופטימיזציה של לוס אחרת. אבל זה לא בדיוק. המאמר כן מציע שיטת אופטימיזציה (מציאת מינימום) לפונקציות קמורות אבל זה בא ממטרה לשפר את Adam או משהו כזה אלא מציע שיטה לשיפור קצב ההתכנסות של אלגוריתם מורד הגרדיאנט (GD) הידוע.  ", 2)

$d.Content.Find.Execute(
    "כנראה שככל יש ברשת יותר מומחים בעלי אותה הארכיטקטורה וכל פעם בוחרים אותו מספר של המומחים הביצועים אמורים להשתפר אולם המחיר הוא המודל גדול יותר.המאמר מנסה לבדוק האם שווה להשתמש בהרבה מאוד במומחים רזים מאוד. המחרים מציעים לעבוד עם מיליון של מומחים של כל אחד מהם היא דל במיוחד. כמובן שכל פעם צריך לבחון את המומחים כל פעם ומכיוון שיש מיליון מומחים אז נדרש מאמץ חישובי לא קטן. המאמר מציע להשתמש בטכניקה הנקראת  product key retrieval כדי להקטין את הסיבוכיות (בגדול זה חלוקה של וקטור המפתחות (keys) לשני חלקים, ביצוע חישוב לכל אחד בנפרד ושילובם).",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "המאמר מתחיל מכך שמבחינה תיאורטית האלגורית של (Polyak-Ruppert (PR הוא זה שאמור להביא התכנסות אופטימלי אבל בפרקטיקה זה פחות קורה (לא ברור לאיזה פרקטיקה הם מתכוונים כי התוצאות שהם נתנו מתייחסות לרשות עמוקות הלא קמורות). PR בעצם עושה אותו GD אבל העדכון האמיתי המוחלק מעריכית עם העדכון האחרון. כלומר באיטרציה t העדכון של GD נכנס עם המקדם 1/t (אפשר לשחק עם זה לפי המאמר אבל קשה להגיע לקצב החלקה אופטימלי).", 2)

$d.Content.Find.Execute(
    "וגיליתי משהו מעניין במאמר הזה - יש scaling law גם ל-MoEs. אולי אסקור אותו בקרוב…",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "המאמר מציע שיטה חדשה (3 שלבים במקום 2 ב-PR) שמשפרת ההתכנסות של PR ללא צורך בבחירה של פרמטר ההחלקה.", 2)

$d.Content.Find.Execute(
    "https://arxiv.org/abs/2407.04153", $true, $false, $false, $false, $false,
    $true, 1, $false, "https://arxiv.org/abs/2405.15682", 2)
